$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A63").Value = "Globo"
$ws.Range("B63").Value = "Inter TV Rural"
$ws.Range("C63").Value = "Agricultura"
$ws.Range("D63").Value = "2025-04-06T00:00"
$ws.Range("E63").Value = "Positivo"
$ws.Range("F63").Value = "Produção de soja em Campos. Agricultores do Norte Fluminense investem no grão tipo exportação. Campos é o maior produtor do Estado do Rio de Janeiro. Uma das vantagens é o baixo frete pela proximidade com o Porto do Açu. Em Santa Cruz, está a maior produção de soja da região. Entrevista com produtor rural José Geraldo Neto; com o secretário de Agricultura, Almy Júnior e com o engenheiro agrônomo, Elias Deulefeu. *matéria* "

$ws.Range("A64").Value = "Globo"
$ws.Range("B64").Value = "RJ TV 2"
$ws.Range("C64").Value = "Defesa Civil"
$ws.Range("D64").Value = "2025-04-05T19:40"
$ws.Range("E64").Value = "Neutro"
$ws.Range("F64").Value = "A chuva no Estado do RJ. Defesas civis de Campos e cidades do Norte Fluminense acompanham deslocamento da frente fria. Repórter *ao vivo*. Chove fraco. A Prefeitura informou que os locais com mais chuva foram Farol e Baixa Grande, na Baixada Campista. Em Dores de Macabu, também choveu bastante. Subsecretário da Defesa Civil, Edison Pessanha, informou que a chuva não causou grandes impactos. "
